$wb = $excel.ActiveWorkbook

# Add a new worksheet for city/subcity masters, placed after the last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "citysubcitymasters"

$ws.Range("A1").Value = "CC0901"
$ws.Range("B1").Value = "Japan"
$ws.Range("C1").Value = "SCC7887"
$ws.Range("D1").Value = "Tokyo"

$ws.Range("A2").Value = "Japan"

$ws.Range("A3").Value = "New Japan"

$ws.Range("A3").Select()
